# Update gh-pages to output generated at 456a3b4
# Applies updated view/fan counts (and one new-event replacement) across
# the four worksheets of the 杭州-漫展信息 workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: 展览 (Exhibitions)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F5").Value  = 349
$ws1.Range("F6").Value  = 544
$ws1.Range("F8").Value  = 11545
$ws1.Range("F12").Value = 2087
$ws1.Range("F19").Value = 156
$ws1.Range("F21").Value = 727
$ws1.Range("F22").Value = 610
$ws1.Range("F25").Value = 714
$ws1.Range("F26").Value = 3425
$ws1.Range("F28").Value = 778
$ws1.Range("F31").Value = 20
$ws1.Range("F32").Value = 967
$ws1.Range("F34").Value = 55
$ws1.Range("F37").Value = 15

# Row 38 - the previous event was replaced by a new one
$ws1.Range("C38").Value = "杭州·百鬼行代号鸢同人only"
$ws1.Range("D38").Value = "保淑路2号 The Queen皇后"
$ws1.Range("E38").Value = "2024.11.02 12:30-11.02 18:40"
$ws1.Range("F38").Value = 1
$ws1.Range("G38").Value = 140
$ws1.Range("H38").Value = "https://show.bilibili.com/platform/detail.html?id=92790"
$ws1.Range("I38").Value = "//i2.hdslb.com/bfs/openplatform/202409/bsDHN4VK1726910606937.jpeg"

$ws1.Range("F39").Value = 4357
$ws1.Range("F40").Value = 5440
$ws1.Range("F42").Value = 112
$ws1.Range("F43").Value = 28
$ws1.Range("F44").Value = 150
$ws1.Range("F47").Value = 23
$ws1.Range("F48").Value = 4089
$ws1.Range("F49").Value = 92

# ---------------------------------------------------------------
# Sheet: 演出 (Performances)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Range("F3").Value  = 4139
$ws2.Range("F11").Value = 576

# ---------------------------------------------------------------
# Sheet: 本地生活 (Local Life)
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")

$ws3.Range("F2").Value = 748
$ws3.Range("F4").Value = 56

# ---------------------------------------------------------------
# Sheet: 全部类型 (All Types)
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value  = 748
$ws4.Range("F4").Value  = 56
$ws4.Range("F7").Value  = 349
$ws4.Range("F8").Value  = 544
$ws4.Range("F9").Value  = 11545
$ws4.Range("F12").Value = 2087
$ws4.Range("F18").Value = 156
$ws4.Range("F20").Value = 4139
$ws4.Range("F22").Value = 727
$ws4.Range("F23").Value = 714
$ws4.Range("F25").Value = 778
$ws4.Range("F29").Value = 20
$ws4.Range("F30").Value = 967
$ws4.Range("F31").Value = 55
$ws4.Range("F34").Value = 15
$ws4.Range("F35").Value = 4357
$ws4.Range("F38").Value = 112
$ws4.Range("F39").Value = 28
$ws4.Range("F40").Value = 150
$ws4.Range("F45").Value = 4089
$ws4.Range("F48").Value = 92
